$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (columns B:E)
$ws.Range("B2").Value = 18.120275670679284
$ws.Range("C2").Value = 4.9842303083526174
$ws.Range("D2").Value = 4.2774990998336033
$ws.Range("E2").Value = 0.68550585053288438

# Row 3 data values (columns B:E)
$ws.Range("B3").Value = 31.578220604750829
$ws.Range("C3").Value = 4.1911312703011223
$ws.Range("D3").Value = -5.755704824788956
$ws.Range("E3").Value = 5.3759646401996122

# Update the selected range to match the new selection highlighted in the workbook
$ws.Range("B1:E3").Select()
